$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1, matching style of existing headers (B1:E1)
$ws.Range("F1").Value = "OSMO_DEF"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.LineStyle = 1

# Fill F2:F4 with "[]" value (plain, unstyled like C/E columns)
$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"
$ws.Range("F4").Value = "[]"
